# Splits the text of a paragraph's single run into several sibling runs,
# one per whitespace-delimited "word" (the separating spaces become their
# own runs too), while leaving the visible text/formatting of the
# paragraph completely unchanged.
#
# A plain Range.Text re-assignment (or InsertBefore/InsertAfter) gets
# silently re-merged by Word's run-coalescing logic whenever the
# resulting runs carry identical formatting, so it can't produce a
# lasting run boundary. Re-assigning a sub-range's FormattedText to
# itself (Range.FormattedText = Range.FormattedText), however, forces
# the engine to record that sub-range as an explicit, separate run
# without introducing any direct character formatting (no stray
# <w:rPr/> survives), because the "new" formatting written is exactly
# the formatting already in effect (inherited from the paragraph style)
# rather than a literal property toggle.
function Split-RangeHere($doc, $absoluteOffset) {
    $r = $doc.Range(0, $absoluteOffset)
    $ft = $r.FormattedText
    $r.FormattedText = $ft
}

function Split-ParagraphIntoWordRuns($doc, $paragraphIndex) {
    $para = $doc.Paragraphs($paragraphIndex)
    $paraRange = $para.Range
    $fullText = $paraRange.Text

    # Paragraph.Range.Text includes the trailing paragraph mark; strip it
    # off before splitting into words.
    $mark = [char]13
    $text = $fullText.TrimEnd($mark)

    $words = $text.Split(" ")

    $start = $paraRange.Start
    $offset = 0
    for ($i = 0; $i -lt $words.Length; $i++) {
        $offset += $words[$i].Length
        if ($i -lt $words.Length - 1) {
            # run boundary right after this word (before the separating space)
            Split-RangeHere $doc ($start + $offset)

            # run boundary right after the separating space itself
            $offset += 1
            Split-RangeHere $doc ($start + $offset)
        }
    }
}

$d = $word.ActiveDocument

Split-ParagraphIntoWordRuns $d 1   # Title: "Answers: Hypothesis Tests"
Split-ParagraphIntoWordRuns $d 2   # Author: "Ellie Trace"
Split-ParagraphIntoWordRuns $d 4   # Abstract: "Answers to questions relating to the guide on Hypothesis tests."
